# Apply the "Updated cryptos list" refresh (Thu Dec 7 18:50:59 UTC 2023).
# Only the B/C/D/E cells that actually changed are touched; D-column price
# cells are forced to stay plain text (matching the original inlineStr cells:
# e.g. "1.50"/"43.399.72" must not be reinterpreted as numbers and lose
# their formatting), then the cell style is reset to "Normal" so no stray
# quote-prefix / number-format style sticks to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D2").Value = "43.399.72"
$ws.Range("D2").Style = "Normal"     # drop the temporary text style again
$ws.Range("E2").Value = "  -1.51%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D3").Value = "2.352.69"
$ws.Range("D3").Style = "Normal"     # drop the temporary text style again
$ws.Range("E3").Value = "  +3.38%  "

# Row 4
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D5").Value = "232.19"
$ws.Range("D5").Style = "Normal"     # drop the temporary text style again
$ws.Range("E5").Value = "  +0.49%  "

# Row 6
$ws.Range("E6").Value = "  +2.35%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D7").Value = "65.98"
$ws.Range("D7").Style = "Normal"     # drop the temporary text style again
$ws.Range("E7").Value = "  +3.64%  "

# Row 8
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D9").Value = "0.452"
$ws.Range("D9").Style = "Normal"     # drop the temporary text style again
$ws.Range("E9").Value = "  +0.68%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D10").Value = "0.0959"
$ws.Range("D10").Style = "Normal"     # drop the temporary text style again
$ws.Range("E10").Value = "  -3.57%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D11").Value = "56.95"
$ws.Range("D11").Style = "Normal"     # drop the temporary text style again
$ws.Range("E11").Value = "  -1.09%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D12").Value = "26.81"
$ws.Range("D12").Style = "Normal"     # drop the temporary text style again
$ws.Range("E12").Value = "  -1.60%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D13").Value = "2.691.62"
$ws.Range("D13").Style = "Normal"     # drop the temporary text style again
$ws.Range("E13").Value = "  +2.99%  "

# Row 14
$ws.Range("E14").Value = "  -1.30%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D15").Value = "15.42"
$ws.Range("D15").Style = "Normal"     # drop the temporary text style again

# Row 17
$ws.Range("D17").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D17").Value = "0.849"
$ws.Range("D17").Style = "Normal"     # drop the temporary text style again
$ws.Range("E17").Value = "  +1.27%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D18").Value = "2.344.75"
$ws.Range("D18").Style = "Normal"     # drop the temporary text style again
$ws.Range("E18").Value = "  +3.83%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D19").Value = "43.337.71"
$ws.Range("D19").Style = "Normal"     # drop the temporary text style again
$ws.Range("E19").Value = "  -1.34%  "

# Row 20
$ws.Range("E20").Value = "  -2.16%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D21").Value = "74.17"
$ws.Range("D21").Style = "Normal"     # drop the temporary text style again
$ws.Range("E21").Value = "  +0.47%  "

# Row 22
$ws.Range("E22").Value = "  +1.84%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D23").Value = "249.39"
$ws.Range("D23").Style = "Normal"     # drop the temporary text style again
$ws.Range("E23").Value = "  -1.12%  "

# Row 24
$ws.Range("E24").Value = "  +16.60%  "

# Row 25
$ws.Range("E25").Value = "  -0.03%  "

# Row 26
$ws.Range("E26").Value = "  -0.67%  "

# Row 27
$ws.Range("E27").Value = "  +0.81%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D28").Value = "9.93"
$ws.Range("D28").Style = "Normal"     # drop the temporary text style again
$ws.Range("E28").Value = "  -1.44%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D29").Value = "174.96"
$ws.Range("D29").Style = "Normal"     # drop the temporary text style again
$ws.Range("E29").Value = "  +1.83%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D30").Value = "22.21"
$ws.Range("D30").Style = "Normal"     # drop the temporary text style again
$ws.Range("E30").Value = "  +6.06%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D31").Value = "1.50"
$ws.Range("D31").Style = "Normal"     # drop the temporary text style again
$ws.Range("E31").Value = "  +7.17%  "

# Row 32
$ws.Range("E32").Value = "  -7.19%  "

# Row 33
$ws.Range("E33").Value = "  +0.64%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D34").Value = "4.99"
$ws.Range("D34").Style = "Normal"     # drop the temporary text style again
$ws.Range("E34").Value = "  +3.86%  "

# Row 35
$ws.Range("E35").Value = "  -2.39%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D36").Value = "4.97"
$ws.Range("D36").Style = "Normal"     # drop the temporary text style again
$ws.Range("E36").Value = "  +1.96%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D37").Value = "2.55"
$ws.Range("D37").Style = "Normal"     # drop the temporary text style again
$ws.Range("E37").Value = "  +9.78%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D38").Value = "6.46"
$ws.Range("D38").Style = "Normal"     # drop the temporary text style again
$ws.Range("E38").Value = "  -1.18%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D39").Value = "3.61"
$ws.Range("D39").Style = "Normal"     # drop the temporary text style again
$ws.Range("E39").Value = "  -5.20%  "

# Row 40
$ws.Range("E40").Value = "  -3.04%  "

# Row 41
$ws.Range("E41").Value = "  +8.85%  "

# Row 42
$ws.Range("E42").Value = "  +0.03%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D43").Value = "18.11"
$ws.Range("D43").Style = "Normal"     # drop the temporary text style again
$ws.Range("E43").Value = "  +2.55%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D44").Value = "1.18"
$ws.Range("D44").Style = "Normal"     # drop the temporary text style again
$ws.Range("E44").Value = "  +8.98%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D45").Value = "99.11"
$ws.Range("D45").Style = "Normal"     # drop the temporary text style again

# Row 46
$ws.Range("E46").Value = "  -0.87%  "

# Row 47: Cronos -> FTXToken
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D47").Value = "4.40"
$ws.Range("D47").Style = "Normal"     # drop the temporary text style again
$ws.Range("E47").Value = "  +0.67%  "

# Row 48: FTXToken -> Cronos
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D48").Value = "0.0948"
$ws.Range("D48").Style = "Normal"     # drop the temporary text style again
$ws.Range("E48").Value = "  -3.97%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D49").Value = "1.439.78"
$ws.Range("D49").Style = "Normal"     # drop the temporary text style again
$ws.Range("E49").Value = "  -0.54%  "

# Row 50
$ws.Range("E50").Value = "  -11.17%  "

# Row 51: Celestia -> RocketPoolETH
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D51").Value = "2.566.45"
$ws.Range("D51").Style = "Normal"     # drop the temporary text style again
$ws.Range("E51").Value = "  +3.18%  "
